$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Re-order / re-caption the nine team-member lines in the title-page table.
#    These are document paragraphs 9..17 (1-based), each a single run with a
#    single <w:t>. We overwrite the paragraph Range.Text directly (not via
#    Find/Replace) because several of the new captions re-use a name that is
#    also one of the *old* search strings a few rows away, which would make a
#    sequential Find/Replace corrupt earlier results.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs

$paras.Item(9).Range.Text  = "Mateusz Ciechan – Product owner"
$paras.Item(10).Range.Text = "Kamil Oleszek – Scrum master"
$paras.Item(11).Range.Text = "Łukasz Ujma - Architekt"
$paras.Item(12).Range.Text = "Aleksander Giera – Architekt"
$paras.Item(13).Range.Text = "Alicja Frankowicz – Programista/Grafik"
$paras.Item(14).Range.Text = "Grzegorz Ryniak – Programista/Grafik"
$paras.Item(15).Range.Text = "Karol Cwynar - Programista"
$paras.Item(16).Range.Text = "Elżbieta Dziedzic - Tester"
$paras.Item(17).Range.Text = "Michał Niemczyk - Tester"

# ---------------------------------------------------------------------------
# 2. The stray "_GoBack" bookmark that used to sit in the TOC is removed from
#    there and re-created as an empty paragraph appended after the list of
#    names (last row of the title table).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$paras = $d.Paragraphs
$lastNamePara = $paras.Item(17)
$anchorEnd = $lastNamePara.Range.End - 1
$lastNamePara.Range.InsertAfter([char]13)

$paras = $d.Paragraphs
$newPara = $paras.Item(18)
$bookmarkRange = $d.Range($anchorEnd, $newPara.Range.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ---------------------------------------------------------------------------
# 3. The footer's cached PAGE field result changes from "5" to "2".
# ---------------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("5", $true, $false, $false, $false, $false, $true,
                            1, $false, "2", 2) | Out-Null
